{"js": "// 1) Update the letter date: \"September 19, 2025\" -> \"September 21, 2025\"\nconst dateResults = context.document.body.search(\"September 19, 2025\", { matchCase: true });\ndateResults.load(\"items\");\nawait context.sync();\nif (dateResults.items.length > 0) {\n  dateResults.items[0].insertText(\"September 21, 2025\", \"Replace\");\n  await context.sync();\n}\n\n// 2) Split the mailing address paragraph \"969 Story Road, San Jose CA 95122\"\n//    into two paragraphs: \"969 Story Road\" and \"San Jose, CA 95122\".\nconst addrResults = context.document.body.search(\"969 Story Road, San Jose CA 95122\", { matchCase: true });\naddrResults.load(\"items\");\nawait context.sync();\n\nif (addrResults.items.length > 0) {\n  // Use the first occurrence, which sits in the address block (the second\n  // occurrence further down is inside the \"PROPERTY ADDRESS:\" line and\n  // must remain unchanged).\n  const addrRange = addrResults.items[0];\n  const addrParagraph = addrRange.paragraphs.getFirst();\n  addrParagraph.load(\"text\");\n  await context.sync();\n\n  // Insert the new \"San Jose, CA 95122\" paragraph right after the address\n  // paragraph, inheriting the same paragraph formatting.\n  addrParagraph.insertParagraph(\"San Jose, CA 95122\", \"After\");\n  await context.sync();\n\n  // Trim the original paragraph's text down to just the street address.\n  addrParagraph.getRange().insertText(\"969 Story Road\", \"Replace\");\n  await context.sync();\n}\n\n// 3) Remove the empty \"No Spacing\" paragraph that immediately follows the\n//    \"...Board of Directors\" signature line.\nconst boardResults = context.document.body.search(\"Board of Directors\", { matchCase: true });\nboardResults.load(\"items\");\nawait context.sync();\n\nif (boardResults.items.length > 0) {\n  const boardParagraph = boardResults.items[0].paragraphs.getFirst();\n  const nextParagraph = boardParagraph.getNext();\n  nextParagraph.load(\"text\");\n  await context.sync();\n\n  if (nextParagraph.text === \"\") {\n    nextParagraph.delete();\n    await context.sync();\n  }\n}\n", "ps1": "$d = $word.ActiveDocument\n\n# --- 1) Update the letter date: \"September 19, 2025\" -> \"September 21, 2025\" ---\n$dateParagraph = $null\n$count = $d.Paragraphs.Count\nfor ($i = 1; $i -le $count; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    $txt = $p.Range.Text.TrimEnd([char]13, [char]7)\n    if ($txt -eq \"September 19, 2025\") {\n        $dateParagraph = $p\n        break\n    }\n}\nif ($dateParagraph -ne $null) {\n    $dateParagraph.Range.Text = \"September 21, 2025\"\n}\n\n# --- 2) Split the mailing address paragraph \"969 Story Road, San Jose CA 95122\"\n#        into two paragraphs: \"969 Story Road\" and \"San Jose, CA 95122\". ---\n$addressParagraph = $null\n$addressIndex = -1\n$count = $d.Paragraphs.Count\nfor ($i = 1; $i -le $count; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    $txt = $p.Range.Text.TrimEnd([char]13, [char]7)\n    if ($txt -eq \"969 Story Road, San Jose CA 95122\") {\n        $addressParagraph = $p\n        $addressIndex = $i\n        break\n    }\n}\nif ($addressParagraph -ne $null) {\n    # Insert a new empty paragraph right after the address paragraph; it\n    # inherits the same paragraph/run formatting (Arial, sz 22) from the\n    # paragraph mark it was split from.\n    $addressParagraph.Range.InsertParagraphAfter()\n\n    $newParagraph = $d.Paragraphs.Item($addressIndex + 1)\n    $newParagraph.Range.Text = \"San Jose, CA 95122\"\n\n    # Trim the original paragraph's text down to just the street address.\n    $addressParagraph.Range.Text = \"969 Story Road\"\n}\n\n# --- 3) Remove the empty \"No Spacing\" paragraph that immediately follows\n#        the \"...Board of Directors\" signature line. ---\n$boardParagraph = $null\n$count = $d.Paragraphs.Count\nfor ($i = 1; $i -le $count; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    if ($p.Range.Text -like \"*Board of Directors*\") {\n        $boardParagraph = $p\n        break\n    }\n}\nif ($boardParagraph -ne $null) {\n    $nextParagraph = $boardParagraph.Next()\n    $nextText = $nextParagraph.Range.Text.TrimEnd([char]13, [char]7)\n    if ($nextText -eq \"\") {\n        $nextParagraph.Range.Delete()\n    }\n}\n"}
